$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "mmWave(InBed)" - append rows 16-24 (In Bed / Occupied)
# ---------------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")
$inBedRows = @(
    @("2026-02-01", "20:04:22", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:23", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:24", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:25", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:26", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:35", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:40", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:44", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:04:45", "20:00", "Bedroom", "In Bed", "Occupied")
)
$startRow = 16
$endRow = $startRow + $inBedRows.Count - 1
$wsInBed.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $inBedRows.Count; $i++) {
    $row = $startRow + $i
    $data = $inBedRows[$i]
    $wsInBed.Cells.Item($row, 1).Value = $data[0]
    $wsInBed.Cells.Item($row, 2).Value = $data[1]
    $wsInBed.Cells.Item($row, 3).Value = $data[2]
    $wsInBed.Cells.Item($row, 4).Value = $data[3]
    $wsInBed.Cells.Item($row, 5).Value = $data[4]
    $wsInBed.Cells.Item($row, 6).Value = $data[5]
}

# ---------------------------------------------------------------------------
# Sheet "mmWave" - append rows 2-3 (Bedroom Door ENTER/EXIT)
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")
$wsMmWave.Range("A2:A3").NumberFormat = "@"

$wsMmWave.Cells.Item(2, 1).Value = "2026-02-01"
$wsMmWave.Cells.Item(2, 2).Value = "20:04:53"
$wsMmWave.Cells.Item(2, 3).Value = "20:00"
$wsMmWave.Cells.Item(2, 4).Value = "Bedroom Door"
$wsMmWave.Cells.Item(2, 5).Value = "ENTER"
$wsMmWave.Cells.Item(2, 6).Value = "User ENTERED Bedroom"

$wsMmWave.Cells.Item(3, 1).Value = "2026-02-01"
$wsMmWave.Cells.Item(3, 2).Value = "20:05:03"
$wsMmWave.Cells.Item(3, 3).Value = "20:00"
$wsMmWave.Cells.Item(3, 4).Value = "Bedroom Door"
$wsMmWave.Cells.Item(3, 5).Value = "EXIT"
$wsMmWave.Cells.Item(3, 6).Value = "User EXITED Bedroom"

# ---------------------------------------------------------------------------
# Sheet "mmWave(BR)" - append rows 13-22 (numeric Value column)
# ---------------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")
$brRows = @(
    @("2026-02-01", "20:04:21", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:04:23", "20:00", "Bedroom", 8, "Occupied"),
    @("2026-02-01", "20:04:24", "20:00", "Bedroom", 3, "Occupied"),
    @("2026-02-01", "20:04:25", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:04:26", "20:00", "Bedroom", 10, "Occupied"),
    @("2026-02-01", "20:04:27", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:04:36", "20:00", "Bedroom", 1, "Occupied"),
    @("2026-02-01", "20:04:40", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:04:44", "20:00", "Bedroom", 22, "Occupied"),
    @("2026-02-01", "20:04:46", "20:00", "Bedroom", 2, "Occupied")
)
$startRow = 13
$endRow = $startRow + $brRows.Count - 1
$wsBR.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $brRows.Count; $i++) {
    $row = $startRow + $i
    $data = $brRows[$i]
    $wsBR.Cells.Item($row, 1).Value = $data[0]
    $wsBR.Cells.Item($row, 2).Value = $data[1]
    $wsBR.Cells.Item($row, 3).Value = $data[2]
    $wsBR.Cells.Item($row, 4).Value = $data[3]
    $wsBR.Cells.Item($row, 5).Value = $data[4]
    $wsBR.Cells.Item($row, 6).Value = $data[5]
}

# ---------------------------------------------------------------------------
# Sheet "mmWave(HR)" - append rows 13-22 (numeric Value column)
# ---------------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")
$hrRows = @(
    @("2026-02-01", "20:04:21", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:04:22", "20:00", "Bedroom", 56, "Occupied"),
    @("2026-02-01", "20:04:23", "20:00", "Bedroom", 51, "Occupied"),
    @("2026-02-01", "20:04:24", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:04:26", "20:00", "Bedroom", 58, "Occupied"),
    @("2026-02-01", "20:04:27", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:04:36", "20:00", "Bedroom", 49, "Occupied"),
    @("2026-02-01", "20:04:40", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:04:44", "20:00", "Bedroom", 70, "Occupied"),
    @("2026-02-01", "20:04:45", "20:00", "Bedroom", 50, "Occupied")
)
$startRow = 13
$endRow = $startRow + $hrRows.Count - 1
$wsHR.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $hrRows.Count; $i++) {
    $row = $startRow + $i
    $data = $hrRows[$i]
    $wsHR.Cells.Item($row, 1).Value = $data[0]
    $wsHR.Cells.Item($row, 2).Value = $data[1]
    $wsHR.Cells.Item($row, 3).Value = $data[2]
    $wsHR.Cells.Item($row, 4).Value = $data[3]
    $wsHR.Cells.Item($row, 5).Value = $data[4]
    $wsHR.Cells.Item($row, 6).Value = $data[5]
}
